$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (row 1 header stays the same text, but we
# rewrite it too for safety/consistency). Clear a generous range so any
# stray formatting/content from rows beyond the new extent is removed.
$ws.Range("A1:I20").Clear()

# Header row (row 1) - bold/style handled by column+row default style
$ws.Range("B1").Value = "Dienstplan"
$ws.Range("C1").Value = "Schicht"
$ws.Range("D1").Value = "Assistenz"
$ws.Range("E1").Value = "Tag"
$ws.Range("F1").Value = "Verfügbarkeit"
$ws.Range("G1").Value = "Dienst"
$ws.Range("H1").Value = "Zuschlag"
$ws.Range("I1").Value = "SchichtTag"

# Row 2 - Dienstplan erstellen
$ws.Range("A2").Value = "Dienstplan erstellen"
$ws.Range("B2").Value = "CRUD"

# Row 3 - Dienstplan generieren
$ws.Range("A3").Value = "Dienstplan generieren"
$ws.Range("B3").Value = "RU"
$ws.Range("F3").Value = "R"
$ws.Range("G3").Value = "CRUD"

# Row 4 - Dienst bearbeiten
$ws.Range("A4").Value = "Dienst bearbeiten"
$ws.Range("B4").Value = "RU"
$ws.Range("F4").Value = "R"
$ws.Range("G4").Value = "RU"

# Row 5 - Dienst tauschen
$ws.Range("A5").Value = "Dienst tauschen"
$ws.Range("B5").Value = "RU"
$ws.Range("F5").Value = "R"
$ws.Range("G5").Value = "RU"

# Row 6 - Schicht erstellen
$ws.Range("A6").Value = "Schicht erstellen"
$ws.Range("C6").Value = "CRUD"
$ws.Range("F6").Value = "R"

# Row 7 - Schicht bearbeiten
$ws.Range("A7").Value = "Schicht bearbeiten"
$ws.Range("C7").Value = "CRUD"
$ws.Range("F7").Value = "R"

# Row 8 - Tag erstellen
$ws.Range("A8").Value = "Tag erstellen"
$ws.Range("E8").Value = "CRUD"

# Row 9 - Tag bearbeiten
$ws.Range("A9").Value = "Tag bearbeiten"
$ws.Range("E9").Value = "CRUD"
$ws.Range("F9").Value = "R"

# Row 10 - Tag zuordnen
$ws.Range("A10").Value = "Tag zuordnen"
$ws.Range("D10").Value = "RU"
$ws.Range("E10").Value = "R"
$ws.Range("F10").Value = "R"
$ws.Range("I10").Value = "CRUD"

# Row 11 - Assistenz erstellen
$ws.Range("A11").Value = "Assistenz erstellen"
$ws.Range("D11").Value = "CRUD"

# Row 12 - Assistenz bearbeiten
$ws.Range("A12").Value = "Assistenz bearbeiten"
$ws.Range("D12").Value = "RUD"

# Row 13 - Assistenz zuordnen
$ws.Range("A13").Value = "Assistenz zuordnen"
$ws.Range("D13").Value = "R"
$ws.Range("F13").Value = "R"
$ws.Range("G13").Value = "CRUD"

# Row 14 - Zuschlag erstellen/ bearbeiten
$ws.Range("A14").Value = "Zuschlag erstellen/ bearbeiten"
$ws.Range("H14").Value = "CRUD"

# Row 15 - Konflikte behandeln
$ws.Range("A15").Value = "Konflikte behandeln"

# Update selected cell to K15 to match the new view state
$ws.Range("K15").Select()
